# Auto-generated data-driven cell updates for paises.xlsx daily refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=1; Col=1; Value="Datos actualizados a 16 de Agosto de 2020 a las 21:17"},
    @{Row=4; Col=2; Value=5546622},
    @{Row=4; Col=3; Value=16833},
    @{Row=4; Col=4; Value=2908631},
    @{Row=4; Col=5; Value=2465121},
    @{Row=4; Col=7; Value=264},
    @{Row=4; Col=8; Value=172870},
    @{Row=5; Col=2; Value=3319845},
    @{Row=5; Col=3; Value=2013},
    @{Row=5; Col=5; Value=808232},
    @{Row=5; Col=7; Value=44},
    @{Row=5; Col=8; Value=107341},
    @{Row=6; Col=2; Value=2647316},
    @{Row=6; Col=3; Value=58108},
    @{Row=6; Col=4; Value=1918076},
    @{Row=6; Col=5; Value=678195},
    @{Row=6; Col=7; Value=961},
    @{Row=6; Col=8; Value=51045},
    @{Row=12; Col=2; Value=385946},
    @{Row=12; Col=3; Value=2044},
    @{Row=12; Col=4; Value=358828},
    @{Row=12; Col=5; Value=16666},
    @{Row=12; Col=7; Value=57},
    @{Row=12; Col=8; Value=10452},
    @{Row=21; Col=2; Value=249309},
    @{Row=21; Col=3; Value=1192},
    @{Row=21; Col=4; Value=230969},
    @{Row=21; Col=5; Value=12366},
    @{Row=21; Col=7; Value=19},
    @{Row=21; Col=8; Value=5974},
    @{Row=22; Col=2; Value=224880},
    @{Row=22; Col=3; Value=402},
    @{Row=22; Col=5; Value=12690},
    @{Row=23; Col=2; Value=218536},
    @{Row=23; Col=3; Value=3015},
    @{Row=23; Col=5; Value=104278},
    @{Row=23; Col=7; Value=1},
    @{Row=23; Col=8; Value=30410},
    @{Row=25; Col=2; Value=161253},
    @{Row=25; Col=3; Value=3420},
    @{Row=25; Col=4; Value=112586},
    @{Row=25; Col=5; Value=46002},
    @{Row=25; Col=7; Value=65},
    @{Row=25; Col=8; Value=2665},
    @{Row=27; Col=2; Value=122039},
    @{Row=27; Col=3; Value=150},
    @{Row=27; Col=4; Value=108436},
    @{Row=27; Col=5; Value=4578},
    @{Row=27; Col=7; Value=1},
    @{Row=27; Col=8; Value=9025},
    @{Row=33; Col=2; Value=92680},
    @{Row=33; Col=3; Value=447},
    @{Row=33; Col=4; Value=68510},
    @{Row=33; Col=5; Value=23485},
    @{Row=33; Col=7; Value=11},
    @{Row=33; Col=8; Value=685},
    @{Row=54; Col=2; Value=42532},
    @{Row=54; Col=3; Value=322},
    @{Row=54; Col=4; Value=40362},
    @{Row=54; Col=5; Value=1939},
    @{Row=55; Col=1; Value="Marruecos"},
    @{Row=55; Col=2; Value=42489},
    @{Row=55; Col=3; Value=1472},
    @{Row=55; Col=4; Value=29344},
    @{Row=55; Col=5; Value=12487},
    @{Row=55; Col=7; Value=26},
    @{Row=55; Col=8; Value=658},
    @{Row=56; Col=1; Value="Kirguistan"},
    @{Row=56; Col=2; Value=41856},
    @{Row=56; Col=3; Value=211},
    @{Row=56; Col=4; Value=34276},
    @{Row=56; Col=5; Value=6085},
    @{Row=56; Col=7; Value=2},
    @{Row=56; Col=8; Value=1495},
    @{Row=57; Col=1; Value="Armenia"},
    @{Row=57; Col=2; Value=41663},
    @{Row=57; Col=3; Value=168},
    @{Row=57; Col=4; Value=34584},
    @{Row=57; Col=5; Value=6261},
    @{Row=57; Col=7; Value=1},
    @{Row=57; Col=8; Value=818},
    @{Row=66; Col=1; Value="Etiopia"},
    @{Row=66; Col=2; Value=29876},
    @{Row=66; Col=3; Value=982},
    @{Row=66; Col=4; Value=12359},
    @{Row=66; Col=5; Value=16989},
    @{Row=66; Col=7; Value=19},
    @{Row=66; Col=8; Value=528},
    @{Row=67; Col=1; Value="Serbia"},
    @{Row=67; Col=2; Value=29682},
    @{Row=67; Col=3; Value=211},
    @{Row=67; Col=4; Value=27061},
    @{Row=67; Col=5; Value=1947},
    @{Row=67; Col=7; Value=4},
    @{Row=67; Col=8; Value=674},
    @{Row=69; Col=2; Value=27257},
    @{Row=69; Col=3; Value=66},
    @{Row=69; Col=5; Value=2119},
    @{Row=91; Col=1; Value="Libano"},
    @{Row=91; Col=2; Value=8881},
    @{Row=91; Col=3; Value=439},
    @{Row=91; Col=4; Value=2724},
    @{Row=91; Col=5; Value=6054},
    @{Row=91; Col=7; Value=6},
    @{Row=91; Col=8; Value=103},
    @{Row=92; Col=1; Value="Guayana Francesa"},
    @{Row=92; Col=2; Value=8588},
    @{Row=92; Col=4; Value=7893},
    @{Row=92; Col=5; Value=642},
    @{Row=92; Col=8; Value=53},
    @{Row=104; Col=2; Value=5785},
    @{Row=104; Col=3; Value=106},
    @{Row=104; Col=4; Value=3349},
    @{Row=104; Col=5; Value=2414},
    @{Row=138; Col=2; Value=1869},
    @{Row=138; Col=3; Value=11},
    @{Row=138; Col=5; Value=326},
    @{Row=138; Col=7; Value=2},
    @{Row=138; Col=8; Value=530},
    @{Row=140; Col=1; Value="Siria"},
    @{Row=140; Col=2; Value=1677},
    @{Row=140; Col=3; Value=84},
    @{Row=140; Col=4; Value=417},
    @{Row=140; Col=5; Value=1196},
    @{Row=140; Col=7; Value=4},
    @{Row=140; Col=8; Value=64},
    @{Row=141; Col=1; Value="Nueva Zelanda"},
    @{Row=141; Col=2; Value=1622},
    @{Row=141; Col=3; Value=13},
    @{Row=141; Col=4; Value=1531},
    @{Row=141; Col=5; Value=69},
    @{Row=141; Col=8; Value=22},
    @{Row=145; Col=1; Value="Republica de Chipre"},
    @{Row=145; Col=2; Value=1339},
    @{Row=145; Col=3; Value=7},
    @{Row=145; Col=4; Value=870},
    @{Row=145; Col=5; Value=449},
    @{Row=145; Col=8; Value=20},
    @{Row=146; Col=1; Value="Georgia"},
    @{Row=146; Col=2; Value=1336},
    @{Row=146; Col=3; Value=15},
    @{Row=146; Col=4; Value=1088},
    @{Row=146; Col=5; Value=231},
    @{Row=146; Col=8; Value=17},
    @{Row=153; Col=2; Value=1166},
    @{Row=153; Col=3; Value=1},
    @{Row=153; Col=5; Value=20},
    @{Row=156; Col=2; Value=1102},
    @{Row=156; Col=3; Value=54},
    @{Row=156; Col=4; Value=200},
    @{Row=156; Col=5; Value=898},
    @{Row=159; Col=2; Value=956},
    @{Row=159; Col=3; Value=4},
    @{Row=159; Col=4; Value=865},
    @{Row=159; Col=5; Value=15},
    @{Row=162; Col=2; Value=855},
    @{Row=162; Col=3; Value=39},
    @{Row=162; Col=5; Value=193},
    @{Row=213; Col=1; Value="Montserrat"},
    @{Row=213; Col=4; Value=12},
    @{Row=213; Col=8; Value=1},
    @{Row=214; Col=1; Value="Islas Malvinas"},
    @{Row=214; Col=4; Value=13},
    @{Row=214; Col=8; Value=0}
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, $u.Col).Value = $u.Value
}

Write-Output "Applied $($updates.Count) cell updates"